$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Sheet1" to "Sheet"
$ws.Name = "Sheet"

# Clear the bold/bordered/centered header style so header cells (A1:D1)
# revert to the default "Normal" style (drops the custom font + border).
$ws.Range("A1:D1").Style = "Normal"

# Set explicit column widths (A=15, B=30, C=507, D=13 "characters").
# ColumnWidth applies a +5/6 rendering offset internally, so subtract it
# here to land on the exact stored width.
$offset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 15 - $offset
$ws.Columns.Item(2).ColumnWidth = 30 - $offset
$ws.Columns.Item(3).ColumnWidth = 507 - $offset
$ws.Columns.Item(4).ColumnWidth = 13 - $offset

# Clear the "Appreciated" (D) column contents for every data row (2-32).
$ws.Range("D2:D32").ClearContents()
